$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rest Assured")

# Row 2: expected response time/size assertions
$ws.Range("F2").Value = 10
$ws.Range("G2").Value = 6

# Row 3: fix typo'd email in the request payload JSON body, bump expected response time
$ws.Range("C3").Value = @"
{
    "first_name": "Ali",
    "last_name": "Ahmad",
    "email": "ali.ahmad2131@gmail.com",
    "password": "12345",
    "confirm_password": "12345"
}
"@
$ws.Range("F3").Value = 10

# Row 4: bump expected response time
$ws.Range("F4").Value = 10
